$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 16.322
$ws.Range("D3").Value = -7.834000000000001
$ws.Range("E3").Value = 16.838
$ws.Range("D4").Value = -7.769
$ws.Range("E9").Value = 17.145
$ws.Range("C11").Value = -12.27
$ws.Range("C12").Value = -10.771
$ws.Range("D14").Value = -7.597999999999999
$ws.Range("C15").Value = -13.023
$ws.Range("E15").Value = 16.182
$ws.Range("E19").Value = 16.35
$ws.Range("E20").Value = 16.31
$ws.Range("E25").Value = 17.135
$ws.Range("D26").Value = -7.825000000000001
$ws.Range("C27").Value = -13.151
$ws.Range("E27").Value = 16.517
$ws.Range("C28").Value = -13.082
$ws.Range("E28").Value = 16.634
$ws.Range("E30").Value = 16.354
$ws.Range("C31").Value = -12.857
$ws.Range("D31").Value = -7.784999999999999
$ws.Range("C32").Value = -13.244
$ws.Range("E32").Value = 16.533
$ws.Range("D35").Value = -7.854000000000001
$ws.Range("C36").Value = -12.721
$ws.Range("D37").Value = -7.632
$ws.Range("C38").Value = -12.608
$ws.Range("D39").Value = -7.218000000000001
$ws.Range("D40").Value = -7.853999999999999
$ws.Range("E44").Value = 16.489
$ws.Range("D45").Value = -7.94
$ws.Range("C46").Value = -13.817
$ws.Range("E47").Value = 16.386
$ws.Range("D52").Value = -7.385000000000001
$ws.Range("C54").Value = -13.141
$ws.Range("C55").Value = -13.53
$ws.Range("C56").Value = -13.364
$ws.Range("D57").Value = -8.451000000000001
$ws.Range("E58").Value = 16.456
$ws.Range("E62").Value = 16.251
$ws.Range("C67").Value = -11.661
$ws.Range("C69").Value = -10.751
$ws.Range("C72").Value = -11.555
$ws.Range("C73").Value = -12.464
$ws.Range("E77").Value = 16.598
$ws.Range("E78").Value = 16.391
$ws.Range("D81").Value = -7.085000000000001
$ws.Range("C83").Value = -13.314
$ws.Range("D83").Value = -8.413999999999998
$ws.Range("E84").Value = 16.438
$ws.Range("C86").Value = -13.846
$ws.Range("E89").Value = 17.142
$ws.Range("C91").Value = -11.175
$ws.Range("E91").Value = 16.983
$ws.Range("E92").Value = 16.883
$ws.Range("C93").Value = -11.979
$ws.Range("E96").Value = 16.566
$ws.Range("C99").Value = -12.635
$ws.Range("D100").Value = -7.998
$ws.Range("D102").Value = -7.642
$ws.Range("E102").Value = 16.321
